# Update picture "alt text" (descr attribute) from the old CloudFront URLs
# to the new local "images/..." paths, per the gh-pages build commit.
#
#   slide 2  : Eso_Anatomy_Labels.png
#   slide 16 : gast_total.png
#   slide 17 : gast_proximal_tumor.png  AND  gast_total.png
#
# NOTE: slide 17 has two <p:pic> shapes that both carry id="0" in the
# source OOXML (a pre-existing authoring quirk). The host resolves shapes
# by that id, so plain `Shapes.Item(2)`/`Item(3)` property reads/writes
# both collide on the first "Picture 1" (Shapes.Item() for *reading*
# properties always lands on the first shape sharing that id). To reach
# the second picture independently we duplicate the colliding shape
# (duplicates get a fresh unique id), delete the original, fix the
# duplicate's position back to the original's, restore the original
# z-order, and then set each AlternativeText independently.

$p = $ppt.ActivePresentation

# ---- Slide 2: single picture, direct update -------------------------
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.Type -eq 13) {
        $sh.AlternativeText = "images/Eso_Anatomy_Labels.png"
    }
}

# ---- Slide 16: single picture, direct update -------------------------
$s16 = $p.Slides.Item(16)
for ($i = 1; $i -le $s16.Shapes.Count; $i++) {
    $sh = $s16.Shapes.Item($i)
    if ($sh.Type -eq 13) {
        $sh.AlternativeText = "images/gast_total.png"
    }
}

# ---- Slide 17: two pictures sharing a colliding shape id --------------
$s17 = $p.Slides.Item(17)

# Find index of the first picture shape (currently resolves to the
# "gast_proximal_tumor.png" picture, which is the one reachable via the
# object model because of the id collision).
$firstPicIndex = -1
for ($i = 1; $i -le $s17.Shapes.Count; $i++) {
    if ($s17.Shapes.Item($i).Type -eq 13) {
        $firstPicIndex = $i
        break
    }
}

$firstPic = $s17.Shapes.Item($firstPicIndex)

# Remember its geometry so the re-created shape can be put back exactly
# where it was.
$origLeft = $firstPic.Left
$origTop = $firstPic.Top
$origWidth = $firstPic.Width
$origHeight = $firstPic.Height

# Duplicate it -- the duplicate gets its own unique shape id, so it (and
# the original left behind) become independently addressable.
$dup = $firstPic.Duplicate()

# Delete the original (reachable only through the colliding id) - this
# leaves the *other* original picture ("gast_total.png") as the sole
# shape with that id, making it independently addressable too.
$s17.Shapes.Item($firstPicIndex).Delete()

# Restore the duplicate's original position/size (Duplicate() offsets it).
$dup.Left = $origLeft
$dup.Top = $origTop
$dup.Width = $origWidth
$dup.Height = $origHeight

# Put it back to its original stacking position (right after the title,
# before the remaining picture) so shape order matches the source file.
$dup.ZOrder(3)  # msoSendBackward

# Now set alt text independently on each picture.
$dup.AlternativeText = "images/gast_proximal_tumor.png"
for ($i = 1; $i -le $s17.Shapes.Count; $i++) {
    $sh = $s17.Shapes.Item($i)
    if ($sh.Type -eq 13 -and $sh.Id -ne $dup.Id) {
        $sh.AlternativeText = "images/gast_total.png"
    }
}
